$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 996.5
$ws.Range("I32").Value = 999
$ws.Range("J32").Value = 996
$ws.Range("K32").Value = 999
$ws.Range("L32").Value = 996
$ws.Range("M32").Value = -673
$ws.Range("N32").Value = -1648
# Row 70
$ws.Range("H70").Value = 940.8387
$ws.Range("I70").Value = 891.625
$ws.Range("J70").Value = 993.3333
$ws.Range("K70").Value = 2674.875
$ws.Range("L70").Value = 2979.9999
$ws.Range("M70").Value = -2404.875
$ws.Range("N70").Value = -3519.9999
# Row 73
$ws.Range("H73").Value = 940.8387
$ws.Range("I73").Value = 891.625
$ws.Range("J73").Value = 993.3333
$ws.Range("K73").Value = 2674.875
$ws.Range("L73").Value = 2979.9999
$ws.Range("M73").Value = -1738.875
$ws.Range("N73").Value = -4851.9999
# Row 80
$ws.Range("H80").Value = 2214.85
$ws.Range("I80").Value = 559.8461
$ws.Range("J80").Value = 5288.4287
$ws.Range("K80").Value = 1679.5383
$ws.Range("L80").Value = 15865.2861
$ws.Range("M80").Value = -681.5382999999999
$ws.Range("N80").Value = -17861.2861
# Row 83
$ws.Range("H83").Value = 2214.85
$ws.Range("I83").Value = 559.8461
$ws.Range("J83").Value = 5288.4287
$ws.Range("K83").Value = 5038.6149
$ws.Range("L83").Value = 47595.85830000001
$ws.Range("M83").Value = -46.61489999999958
$ws.Range("N83").Value = -57579.85830000001
# Row 137
$ws.Range("H137").Value = 1745.6666
$ws.Range("I137").Value = 1001.25
$ws.Range("J137").Value = 2341.2
$ws.Range("K137").Value = 3003.75
$ws.Range("L137").Value = 7023.599999999999
$ws.Range("M137").Value = -453.75
$ws.Range("N137").Value = -12123.6
# Row 138
$ws.Range("H138").Value = 2229.2856
$ws.Range("I138").Value = 1787.72
$ws.Range("J138").Value = 2585.3872
$ws.Range("K138").Value = 5363.16
$ws.Range("L138").Value = 7756.1616
$ws.Range("M138").Value = -223.1599999999999
$ws.Range("N138").Value = -18036.1616

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 353648.88
$ws.Range("I32").Value = 410273.9
$ws.Range("J32").Value = 13898.667
$ws.Range("K32").Value = 410273.9
$ws.Range("L32").Value = 13898.667
$ws.Range("M32").Value = -409986.9
$ws.Range("N32").Value = -14472.667
# Row 45
$ws.Range("H45").Value = 3447.8823
$ws.Range("I45").Value = 2532
$ws.Range("K45").Value = 2532
$ws.Range("M45").Value = -2155
# Row 110
$ws.Range("H110").Value = 2820.125
$ws.Range("I110").Value = 3000.6924
$ws.Range("K110").Value = 3000.6924
$ws.Range("M110").Value = -955.6923999999999
# Row 122
$ws.Range("H122").Value = 1327.7894
$ws.Range("I122").Value = 1160
$ws.Range("J122").Value = 1957
$ws.Range("K122").Value = 3480
$ws.Range("L122").Value = 5871
$ws.Range("M122").Value = -1030
$ws.Range("N122").Value = -10771

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 733.5
$ws.Range("I16").Value = 680.2
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 680.2
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -393.2
$ws.Range("N16").Value = -1574
# Row 107
$ws.Range("H107").Value = 2976912
$ws.Range("I107").Value = 5682326.5
$ws.Range("K107").Value = 5682326.5
$ws.Range("M107").Value = -5680406.5
# Row 113
$ws.Range("H113").Value = 733.5
$ws.Range("I113").Value = 680.2
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 680.2
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1489.8
$ws.Range("N113").Value = -5340
# Row 122
$ws.Range("H122").Value = 1996.3684
$ws.Range("I122").Value = 1990.5
$ws.Range("J122").Value = 1997.0588
$ws.Range("K122").Value = 5971.5
$ws.Range("L122").Value = 5991.1764
$ws.Range("M122").Value = -3521.5
$ws.Range("N122").Value = -10891.1764

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 83333390
$ws.Range("I14").Value = 83333390
$ws.Range("K14").Value = 250000170
$ws.Range("M14").Value = -249999997
# Row 68
$ws.Range("H68").Value = 1347.6559
$ws.Range("I68").Value = 824.2105
$ws.Range("J68").Value = 1482.0541
$ws.Range("K68").Value = 2472.6315
$ws.Range("L68").Value = 4446.1623
$ws.Range("M68").Value = -1661.6315
$ws.Range("N68").Value = -6068.1623
# Row 71
$ws.Range("H71").Value = 1347.6559
$ws.Range("I71").Value = 824.2105
$ws.Range("J71").Value = 1482.0541
$ws.Range("K71").Value = 7417.8945
$ws.Range("L71").Value = 13338.4869
$ws.Range("M71").Value = -3361.8945
$ws.Range("N71").Value = -21450.4869
# Row 107
$ws.Range("H107").Value = 1522.7333
$ws.Range("I107").Value = 267.57895
$ws.Range("J107").Value = 3690.7273
$ws.Range("K107").Value = 802.73685
$ws.Range("L107").Value = 11072.1819
$ws.Range("M107").Value = 1117.26315
$ws.Range("N107").Value = -14912.1819
# Row 136
$ws.Range("H136").Value = 3371.9333
$ws.Range("I136").Value = 3175.4443
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 9526.332900000001
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -4426.332900000001
$ws.Range("N136").Value = -21200.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 15332.667
$ws.Range("I43").Value = 6000
$ws.Range("K43").Value = 6000
$ws.Range("M43").Value = -5849
# Row 46
$ws.Range("H46").Value = 5699.45
$ws.Range("J46").Value = 4421
$ws.Range("L46").Value = 4421
$ws.Range("N46").Value = -4733
# Row 80
$ws.Range("H80").Value = 42797490
$ws.Range("I80").Value = 57059220
$ws.Range("J80").Value = 12300
$ws.Range("K80").Value = 57059220
$ws.Range("L80").Value = 12300
$ws.Range("M80").Value = -57058222
$ws.Range("N80").Value = -14296
# Row 83
$ws.Range("H83").Value = 42797490
$ws.Range("I83").Value = 57059220
$ws.Range("J83").Value = 12300
$ws.Range("K83").Value = 285296100
$ws.Range("L83").Value = 61500
$ws.Range("M83").Value = -285291108
$ws.Range("N83").Value = -71484
# Row 102
$ws.Range("H102").Value = 1844
$ws.Range("I102").Value = 1844
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1844
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -222
$ws.Range("N102").ClearContents()
# Row 107
$ws.Range("H107").Value = 1064.1818
$ws.Range("I107").Value = 851.5
$ws.Range("J107").Value = 1185.7142
$ws.Range("K107").Value = 851.5
$ws.Range("L107").Value = 1185.7142
$ws.Range("M107").Value = 1068.5
$ws.Range("N107").Value = -5025.7142

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 35716452
$ws.Range("I7").Value = 47621044
$ws.Range("J7").Value = 2679.2856
$ws.Range("K7").Value = 47621044
$ws.Range("L7").Value = 2679.2856
$ws.Range("M7").Value = -47620932
$ws.Range("N7").Value = -2903.2856
# Row 122
$ws.Range("H122").Value = 4024.652
$ws.Range("I122").Value = 2597.9375
$ws.Range("J122").Value = 7285.7144
$ws.Range("K122").Value = 7793.8125
$ws.Range("L122").Value = 21857.1432
$ws.Range("M122").Value = -5343.8125
$ws.Range("N122").Value = -26757.1432
# Row 126
$ws.Range("H126").Value = 35716452
$ws.Range("I126").Value = 47621044
$ws.Range("J126").Value = 2679.2856
$ws.Range("K126").Value = 142863132
$ws.Range("L126").Value = 8037.8568
$ws.Range("M126").Value = -142860662
$ws.Range("N126").Value = -12977.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1755
$ws.Range("I122").Value = 1453.8462
$ws.Range("J122").Value = 2314.2856
$ws.Range("K122").Value = 4361.5386
$ws.Range("L122").Value = 6942.8568
$ws.Range("M122").Value = -1911.5386
$ws.Range("N122").Value = -11842.8568
